$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSteps")
$wsData = $wb.Worksheets.Item("TestData")

# Insert a new row's worth of content above row 3 ("checkAccessibility" step)
# by shifting the existing rows down one at a time. Work from the bottom up
# so a source row is never overwritten before it has been copied to its new
# destination.
$ws.Range("D8:F8").Copy($ws.Range("D9:F9"))
$ws.Range("D7:F7").Copy($ws.Range("D8:F8"))
$ws.Range("A6:F6").Copy($ws.Range("A7:F7"))
$ws.Range("A5:F5").Copy($ws.Range("A6:F6"))
$ws.Range("A4:F4").Copy($ws.Range("A5:F5"))
$ws.Range("A3:F3").Copy($ws.Range("A4:F4"))

# New row 3 content: an accessibility-check step for this test case.
$ws.Range("A3").Value = "checkAccessibility"
$ws.Range("B3").Value = "TC_SM_FOU_ListView_D1"
$ws.Range("C3").Value = ""

# B3 gets a distinct (no-fill, bordered) style with wrapped text.
$ws.Range("F2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B3").WrapText = $true
$ws.Range("B3").HorizontalAlignment = -4131
$ws.Range("B3").VerticalAlignment = -4160

# Update the remembered selection on each sheet.
$wsData.Range("B13").Select()
$ws.Activate()
$ws.Range("B4").Select()
